$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "63.561.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -1.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.041.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -2.16%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "557.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "141.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  -1.60%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.040.62"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.518"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  +3.39%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.29"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  -12.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.488"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  +5.41%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "35.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.539.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "63.650.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  -1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.047.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -1.94%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.77"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -0.43%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "14.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "14.61"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  +9.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.680"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "7.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "82.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  +1.91%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -0.92%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -2.42%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "26.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  -0.25%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -1.62%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.74"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  -0.47%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "54.68"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.0407"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "443.41"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0811"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "3.009.93"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -0.48%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  +2.24%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "8.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -0.52%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value2 = "TheGraph"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.267"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +2.19%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value2 = "InjectiveProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "27.89"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  +6.55%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  -0.02%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "117.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0₃0512"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  -1.10%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -0.04%  "
